$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, shifting existing rows 130:194 down to 131:195.
$ws.Rows("130:130").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A130").Value2 = 10
$ws.Range("B130").Value2 = "Vega Modelo de Temuco"
$ws.Range("C130").Value2 = "La Araucanía"
$ws.Range("D130").Value2 = 45134
$ws.Range("E130").Value2 = 9
$ws.Range("F130").Value2 = 100112035
$ws.Range("G130").Value2 = "Bruselas (repollito)"
$ws.Range("H130").Value2 = "Sin especificar"
$ws.Range("I130").Value2 = "Primera"
$ws.Range("J130").Value2 = 55
$ws.Range("K130").Value2 = 25000
$ws.Range("L130").Value2 = 25000
$ws.Range("M130").Value2 = 25000
$ws.Range("N130").Value2 = "$/malla 15 kilos"
$ws.Range("O130").Value2 = "Provincia de Quillota"
$ws.Range("P130").Value2 = 1667
$ws.Range("Q130").Value2 = 15
$ws.Range("R130").Value2 = "Hortaliza"
